$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''24.407.26'
$ws.Range('E2').Value = '  +9.27%  '
$ws.Range('D3').Value = '''1.678.48'
$ws.Range('E3').Value = '  +4.96%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '''306.68'
$ws.Range('D6').Value = '''0.9969'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '''0.3707'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '''0.3439'
$ws.Range('E8').Value = '  +1.29%  '
$ws.Range('D9').Value = '''48.17'
$ws.Range('E9').Value = '  +12.86%  '
$ws.Range('D10').Value = '''1.181'
$ws.Range('E10').Value = '  +3.83%  '
$ws.Range('D11').Value = '''0.07266'
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('D12').Value = '''0.9993'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '''20.42'
$ws.Range('E13').Value = '  +3.76%  '
$ws.Range('D14').Value = '''6.109'
$ws.Range('E14').Value = '  +3.22%  '
$ws.Range('D15').Value = '''6.745'
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = '''1.676.39'
$ws.Range('E16').Value = '  +4.98%  '
$ws.Range('D17').Value = '''0.00001110'
$ws.Range('E17').Value = '  +2.94%  '
$ws.Range('D18').Value = '''0.9971'
$ws.Range('D19').Value = '''0.06718'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('D20').Value = '''81.12'
$ws.Range('E20').Value = '  +3.72%  '
$ws.Range('D21').Value = '''16.43'
$ws.Range('E21').Value = '  +1.79%  '
$ws.Range('D22').Value = '''6.102'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('D23').Value = '''11.96'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').Value = '''24.354.76'
$ws.Range('E24').Value = '  +8.94%  '
$ws.Range('D25').Value = '''2.438'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').Value = '''3.365'
$ws.Range('E26').Value = '  -11.28%  '
$ws.Range('D27').Value = '''2.666'
$ws.Range('E27').Value = '  +6.88%  '
$ws.Range('D28').Value = '''152.26'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('D29').Value = '''19.59'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = '''1.861.21'
$ws.Range('E30').Value = '  +4.86%  '
$ws.Range('D31').Value = '''127.25'
$ws.Range('E31').Value = '  +5.53%  '
$ws.Range('D32').Value = '''6.318'
$ws.Range('E32').Value = '  +4.83%  '
$ws.Range('D33').Value = '''4.027'
$ws.Range('E33').Value = '  -3.17%  '
$ws.Range('D34').Value = '''0.9702'
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('D35').Value = '''1.736'
$ws.Range('E35').Value = '  +8.11%  '
$ws.Range('D36').Value = '''0.08474'
$ws.Range('E36').Value = '  +2.80%  '
$ws.Range('D37').Value = '''9.136'
$ws.Range('E37').Value = '  +6.52%  '
$ws.Range('D38').Value = '''0.06494'
$ws.Range('E38').Value = '  +5.87%  '
$ws.Range('D39').Value = '''12.32'
$ws.Range('E39').Value = '  +4.66%  '
$ws.Range('D40').Value = '''5.342'
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').Value = '''0.02338'
$ws.Range('E41').Value = '  +5.57%  '
$ws.Range('D42').Value = '''1.261'
$ws.Range('E42').Value = '  +2.09%  '
$ws.Range('D43').Value = '''0.2113'
$ws.Range('E43').Value = '  +4.15%  '
$ws.Range('D44').Value = '''0.6185'
$ws.Range('E44').Value = '  +4.86%  '
$ws.Range('D45').Value = '''0.9971'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').Value = '''3.780'
$ws.Range('E46').Value = '  +2.93%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.5951'
$ws.Range('E47').Value = '  +4.39%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''12.97'
$ws.Range('E48').Value = '  -1.73%  '
$ws.Range('D49').Value = '''127.15'
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('D50').Value = '''2.025'
$ws.Range('E50').Value = '  +2.94%  '
$ws.Range('D51').Value = '''0.07215'
$ws.Range('E51').Value = '  +5.65%  '
